# Golf.xlsx - "Lochmere" sheet (xl/worksheets/sheet4.xml)
# Add a new round of scores (week of 2021-04-22) below the existing data,
# and extend the "penalties" row of the prior week with a couple more
# hole-penalty markers (switching to ScoreHoles for many-to-many).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Lochmere")

# ---------------------------------------------------------------------
# New row 30: the date header for the new round (mirrors row 20 / row 25)
# ---------------------------------------------------------------------
$ws.Range("A25").Copy()
$ws.Range("A30").PasteSpecial(-4122)
$ws.Range("A30").Value = 44308

# ---------------------------------------------------------------------
# New row 31: "Blue" tee scores (mirrors row 21 / row 26)
# Only copy formats for the cells that actually carry non-default
# styling in the source row, to avoid materializing stray blank cells.
# ---------------------------------------------------------------------
$ws.Range("A21:C21").Copy()
$ws.Range("A31:C31").PasteSpecial(-4122)
$ws.Range("K21:L21").Copy()
$ws.Range("K31:L31").PasteSpecial(-4122)
$ws.Range("M21").Copy()
$ws.Range("M31").PasteSpecial(-4122)
$ws.Range("U21:V21").Copy()
$ws.Range("U31:V31").PasteSpecial(-4122)
$ws.Range("W21:X21").Copy()
$ws.Range("W31:X31").PasteSpecial(-4122)

$ws.Range("A31").Value = "Blue"
$ws.Range("D31").Value = 8
$ws.Range("E31").Value = 4
$ws.Range("F31").Value = 6
$ws.Range("G31").Value = 5
$ws.Range("H31").Value = 5
$ws.Range("I31").Value = 6
$ws.Range("J31").Value = 6
$ws.Range("K31").Value = 6
$ws.Range("L31").Value = 8
$ws.Range("M31").Value = 54
$ws.Range("N31").Value = 6
$ws.Range("O31").Value = 3
$ws.Range("P31").Value = 6
$ws.Range("Q31").Value = 8
$ws.Range("R31").Value = 10
$ws.Range("S31").Value = 6
$ws.Range("T31").Value = 7
$ws.Range("U31").Value = 6
$ws.Range("V31").Value = 5
$ws.Range("W31").Value = 57
$ws.Range("X31").Value = 111

# ---------------------------------------------------------------------
# New row 32: "putts" row (mirrors row 22 / row 27)
# ---------------------------------------------------------------------
$ws.Range("A22").Copy()
$ws.Range("A32").PasteSpecial(-4122)
$ws.Range("K22:L22").Copy()
$ws.Range("K32:L32").PasteSpecial(-4122)
$ws.Range("M22").Copy()
$ws.Range("M32").PasteSpecial(-4122)
$ws.Range("U22:V22").Copy()
$ws.Range("U32:V32").PasteSpecial(-4122)
$ws.Range("W22:X22").Copy()
$ws.Range("W32:X32").PasteSpecial(-4122)

$ws.Range("A32").Value = "putts"
$ws.Range("D32").Value = 4
$ws.Range("E32").Value = 1
$ws.Range("F32").Value = 3
$ws.Range("G32").Value = 1
$ws.Range("H32").Value = 2
$ws.Range("I32").Value = 1
$ws.Range("J32").Value = 2
$ws.Range("K32").Value = 2
$ws.Range("L32").Value = 2
$ws.Range("M32").Value = 18
$ws.Range("N32").Value = 2
$ws.Range("O32").Value = 3
$ws.Range("P32").Value = 2
$ws.Range("Q32").Value = 2
$ws.Range("R32").Value = 1
$ws.Range("S32").Value = 2
$ws.Range("T32").Value = 1
$ws.Range("U32").Value = 2
$ws.Range("V32").Value = 2
$ws.Range("W32").Value = 17
$ws.Range("X32").Value = 35

# ---------------------------------------------------------------------
# New row 33: "penalties" row (mirrors row 23 / row 28), with a couple
# of extra hole markers (J, Q, R) plus U/V now populated too
# ---------------------------------------------------------------------
$ws.Range("A28").Copy()
$ws.Range("A33").PasteSpecial(-4122)
$ws.Range("K28").Copy()
$ws.Range("K33").PasteSpecial(-4122)
$ws.Range("U28:V28").Copy()
$ws.Range("U33:V33").PasteSpecial(-4122)

$ws.Range("A33").Value = "penalties"
$ws.Range("F33").Value = "W"
$ws.Range("I33").Value = "W"
$ws.Range("J33").Value = "W"
$ws.Range("P33").Value = " "
$ws.Range("Q33").Value = "LB"
$ws.Range("R33").Value = "W"
$ws.Range("T33").Value = "W"
$ws.Range("U33").Value = "W"
$ws.Range("V33").Value = "LB"

# ---------------------------------------------------------------------
# Update view: scroll back to column A and select the new date cell
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("A30").Select()
